# Update "paises.xlsx" (Pais sheet) with the new COVID-19 daily snapshot.
# 1) Push the updated case counters for the countries whose figures moved.
# 2) Refresh the "datos actualizados" timestamp footer in A1.
# 3) Re-sort the country table (rows 4:216) by "Casos totales" (column B) descending,
#    exactly as the source data generator does for every snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Updated country statistics (row numbers are the *current*, pre-sort positions) ---
# Estados Unidos (row 4)
$ws.Cells.Item(4,2).Value  = 767402
$ws.Cells.Item(4,3).Value  = 2766
$ws.Cells.Item(4,4).Value  = 71396
$ws.Cells.Item(4,5).Value  = 654777
$ws.Cells.Item(4,6).Value  = 13566
$ws.Cells.Item(4,7).Value  = 654
$ws.Cells.Item(4,8).Value  = 41229

# Reino Unido (row 9)
$ws.Cells.Item(9,2).Value  = 124743
$ws.Cells.Item(9,3).Value  = 4676
$ws.Cells.Item(9,4).Value  = 0
$ws.Cells.Item(9,5).Value  = 107890
$ws.Cells.Item(9,6).Value  = 1559
$ws.Cells.Item(9,7).Value  = 449
$ws.Cells.Item(9,8).Value  = 16509

# Canada (row 16)
$ws.Cells.Item(16,2).Value = 35662
$ws.Cells.Item(16,3).Value = 606
$ws.Cells.Item(16,4).Value = 12149
$ws.Cells.Item(16,5).Value = 21895
$ws.Cells.Item(16,6).Value = 557
$ws.Cells.Item(16,7).Value = 31
$ws.Cells.Item(16,8).Value = 1618

# Polonia (row 31)
$ws.Cells.Item(31,2).Value = 9593
$ws.Cells.Item(31,3).Value = 306
$ws.Cells.Item(31,4).Value = 1133
$ws.Cells.Item(31,5).Value = 8080
$ws.Cells.Item(31,6).Value = 160
$ws.Cells.Item(31,7).Value = 20
$ws.Cells.Item(31,8).Value = 380

# Chequia (row 39)
$ws.Cells.Item(39,2).Value = 6838
$ws.Cells.Item(39,3).Value = 92
$ws.Cells.Item(39,4).Value = 1559
$ws.Cells.Item(39,5).Value = 5085
$ws.Cells.Item(39,6).Value = 84
$ws.Cells.Item(39,7).Value = 8
$ws.Cells.Item(39,8).Value = 194

# Cuba (row 81)
$ws.Cells.Item(81,2).Value = 1087
$ws.Cells.Item(81,3).Value = 52
$ws.Cells.Item(81,4).Value = 285
$ws.Cells.Item(81,5).Value = 766
$ws.Cells.Item(81,6).Value = 9
$ws.Cells.Item(81,7).Value = 2
$ws.Cells.Item(81,8).Value = 36

# Georgia (row 110)
$ws.Cells.Item(110,2).Value = 402
$ws.Cells.Item(110,3).Value = 8
$ws.Cells.Item(110,4).Value = 95
$ws.Cells.Item(110,5).Value = 303
$ws.Cells.Item(110,6).Value = 6
$ws.Cells.Item(110,7).Value = 0
$ws.Cells.Item(110,8).Value = 4

# Mauricio (row 112)
$ws.Cells.Item(112,2).Value = 328
$ws.Cells.Item(112,3).Value = 0
$ws.Cells.Item(112,4).Value = 224
$ws.Cells.Item(112,5).Value = 95
$ws.Cells.Item(112,6).Value = 3
$ws.Cells.Item(112,7).Value = 0
$ws.Cells.Item(112,8).Value = 9

# Montenegro (row 114)
$ws.Cells.Item(114,2).Value = 312
$ws.Cells.Item(114,3).Value = 4
$ws.Cells.Item(114,4).Value = 88
$ws.Cells.Item(114,5).Value = 219
$ws.Cells.Item(114,6).Value = 7
$ws.Cells.Item(114,7).Value = 0
$ws.Cells.Item(114,8).Value = 5

# Isla de Man (row 116)
$ws.Cells.Item(116,2).Value = 300
$ws.Cells.Item(116,3).Value = 2
$ws.Cells.Item(116,4).Value = 200
$ws.Cells.Item(116,5).Value = 91
$ws.Cells.Item(116,6).Value = 15
$ws.Cells.Item(116,7).Value = 3
$ws.Cells.Item(116,8).Value = 9

# Liberia (row 145)
$ws.Cells.Item(145,2).Value = 99
$ws.Cells.Item(145,3).Value = 8
$ws.Cells.Item(145,4).Value = 7
$ws.Cells.Item(145,5).Value = 84
$ws.Cells.Item(145,6).Value = 0
$ws.Cells.Item(145,7).Value = 0
$ws.Cells.Item(145,8).Value = 8

# --- 2) Refresh the footer timestamp text in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 17:52"

# --- 3) Re-sort the whole country table by "Casos totales" (column B) descending ---
$dataRange = $ws.Range("A4:H216")
$sortKey   = $ws.Range("B4:B216")
$dataRange.Sort($sortKey, 2)
